$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new day column ("14-dec") before column EQ (147) ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Shift everything from EQ onward one column to the right.
$wsPrix.Range("EQ1").EntireColumn.Insert()

# New header cell for the inserted column.
$wsPrix.Cells.Item(1, 147).Value = "14-dec"

# The inserted column has no data yet for this day -> "-" placeholder, like the
# other not-yet-available day columns further to the right.
for ($r = 2; $r -le 25; $r++) {
    $wsPrix.Cells.Item($r, 147).Value = "-"
}

# --- Sheet "Gaz": append the next daily price row ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$rowGaz = 177
$wsGaz.Cells.Item($rowGaz, 1).NumberFormat = "@"
$wsGaz.Cells.Item($rowGaz, 1).Value = "2025-12-12"
$wsGaz.Cells.Item($rowGaz, 1).ClearFormats()
$wsGaz.Cells.Item($rowGaz, 2).Value = 26.195

# --- Sheet "CO2": append the next daily price row ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$rowCo2 = 177
$wsCo2.Cells.Item($rowCo2, 1).NumberFormat = "@"
$wsCo2.Cells.Item($rowCo2, 1).Value = "2025-12-12"
$wsCo2.Cells.Item($rowCo2, 1).ClearFormats()
$wsCo2.Cells.Item($rowCo2, 2).Value = 84.09999999999999
